# Add a "treatment" feature: babies should have a room with the
# specific treatment they need.
#
# This adds a new "treatment" header column to both the "babies" sheet
# (new column D) and the "rooms" sheet (new column I), and moves the
# active selection/tab to reflect where the new header was typed
# (babies!D2, with the babies tab now active; rooms!I1 remembers its
# own last selection on the new header cell).

$wb = $excel.ActiveWorkbook

$wsBabies = $wb.Worksheets.Item("babies")
$wsRooms  = $wb.Worksheets.Item("rooms")

# New "treatment" header on babies (column D, first empty column after C)
$wsBabies.Range("D1").Value = "treatment"

# New "treatment" header on rooms (column I, first empty column after H)
$wsRooms.Range("I1").Value = "treatment"

# Leave the rooms sheet selection parked on the new header cell.
$wsRooms.Activate()
$wsRooms.Range("I1").Select()

# Babies becomes the active tab, with the selection sitting just below
# the new header, ready for data entry.
$wsBabies.Activate()
$wsBabies.Range("D2").Select()
